# Apply the "contingencies with rene fine" edit:
#   - B1 = 0              (bold, thin-bordered, centered/top-aligned cell)
#   - A2 = 0              (same formatting as B1)
#   - B2 = "disconnected_elements"  (plain text -> becomes a shared string)
#
# This grows the sheet's dimension to A1:B2, adds a bold font + thin border
# cell style (fontId=1 / borderId=1 / cellXfs index 1) and a sharedStrings
# table containing "disconnected_elements".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- values -----------------------------------------------------------
$ws.Range("B1").Value = 0
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "disconnected_elements"

# --- formatting for B1 --------------------------------------------------
$r1 = $ws.Range("B1")
$r1.Font.Bold = $true
$r1.HorizontalAlignment = -4108   # xlCenter
$r1.VerticalAlignment   = -4160   # xlTop
$r1.Borders.LineStyle = 1         # xlContinuous
$r1.Borders.Weight    = 2         # xlThin

# --- give A2 the exact same style as B1 --------------------------------
# (copy/paste-special-formats keeps the style table minimal: applying the
# same sequence of property writes to A2 independently creates a second,
# slightly different cellXfs entry in this runtime, so we clone B1's
# formatting instead, which re-uses the very same style index.)
$r1.Copy()
$ws.Range("A2").PasteSpecial(-4122)  # xlPasteFormats
